$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A scratch cell used to "type" text that looks like a date (e.g. "02-10-2021")
# without Excel's automatic date recognition converting it to a serial date.
# Entering it as a formula that evaluates to a text string, then copying and
# pasting only the resulting VALUE into the real cell, preserves it as plain
# text (a shared string) and keeps cell formatting/styles untouched.
$scratch = $ws.Cells.Item(1000, 1)

function Set-TextValue($cell, $text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 275 (01-10-2021) previously had no BCP/BCU figures recorded - fill
# them in now, matching every other row's values.
$ws.Cells.Item(275, 2).Value = 187
$ws.Cells.Item(275, 3).Value = 628

# Append the new daily rows for 02-10-2021 .. 06-10-2021.
$newRows = @(
    @{ Date = "02-10-2021"; HasBcpBcu = $true  },
    @{ Date = "03-10-2021"; HasBcpBcu = $true  },
    @{ Date = "04-10-2021"; HasBcpBcu = $true  },
    @{ Date = "05-10-2021"; HasBcpBcu = $true  },
    @{ Date = "06-10-2021"; HasBcpBcu = $false }
)

$row = 276
foreach ($r in $newRows) {
    Set-TextValue $ws.Cells.Item($row, 1) $r.Date
    if ($r.HasBcpBcu) {
        $ws.Cells.Item($row, 2).Value = 187
        $ws.Cells.Item($row, 3).Value = 628
    }
    $ws.Cells.Item($row, 4).Value = 3940
    $ws.Cells.Item($row, 5).Value = 30
    $row++
}

# Clean up the scratch cell so it leaves no trace in the saved workbook.
$scratch.Clear()
